$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/11/2024  Through  11/17/2024"

# --- Cells changing between numeric and text ("0" / "***.*") representation ---
# Set value first (text values get a leading apostrophe so they are not
# auto-coerced back to numbers), then copy number-format/font from a
# neighbouring cell that already has the desired style so the resulting
# style index matches the sibling cells using that style.
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C17").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("G23").Value = "'0"
$ws.Range("G14").Copy()
$ws.Range("G23").PasteSpecial(-4122)

$ws.Range("H23").Value = "'***.*"
$ws.Range("H14").Copy()
$ws.Range("H23").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D28").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 100
$ws.Range("E16").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -23.893805309734
$ws.Range("L16").Value = -32.283464566929
$ws.Range("M16").Value = -27.118644067796
$ws.Range("N16").Value = -84.162062615101
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 66.666666666666
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 17.821782178217
$ws.Range("L17").Value = -17.361111111111
$ws.Range("M17").Value = 77.611940298507
$ws.Range("N17").Value = -49.576271186440
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 155
$ws.Range("K18").Value = -14.193548387096
$ws.Range("L18").Value = -43.644067796610
$ws.Range("M18").Value = -37.558685446009
$ws.Range("N18").Value = -86.973555337904
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 589
$ws.Range("J19").Value = 609
$ws.Range("K19").Value = -3.284072249589
$ws.Range("L19").Value = 8.872458410351
$ws.Range("M19").Value = 122.264150943396
$ws.Range("N19").Value = 79.027355623100
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 12
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 105
$ws.Range("J20").Value = 143
$ws.Range("K20").Value = -26.573426573426
$ws.Range("L20").Value = -36.746987951807
$ws.Range("M20").Value = -20.454545454545
$ws.Range("N20").Value = -86.924034869240
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 107
$ws.Range("H21").Value = -3.738317757009
$ws.Range("I21").Value = 1042
$ws.Range("J21").Value = 1127
$ws.Range("K21").Value = -7.542147293700
$ws.Range("L21").Value = -15.077424612876
$ws.Range("M21").Value = 30.904522613065
$ws.Range("N21").Value = -64.642008822531
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 0
$ws.Range("F23").Value = 2
$ws.Range("L23").Value = -10
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 23.188405797101
$ws.Range("I24").Value = 877
$ws.Range("J24").Value = 798
$ws.Range("K24").Value = 9.899749373433
$ws.Range("L24").Value = -3.837719298245
$ws.Range("M24").Value = 60.036496350365
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = 17.647058823529
$ws.Range("I25").Value = 551
$ws.Range("J25").Value = 442
$ws.Range("K25").Value = 24.660633484162
$ws.Range("L25").Value = 2.990654205607
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 45
$ws.Range("I26").Value = 248
$ws.Range("J26").Value = 217
$ws.Range("K26").Value = 14.285714285714
$ws.Range("L26").Value = 0.813008130081
$ws.Range("M26").Value = 31.216931216931
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 39
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -7.142857142857
$ws.Range("L28").Value = 44.444444444444
